# Auto-generated Excel COM-interop script to apply scheduled-runner data refresh
# Updates currentAveragePrice / Leve price / profit columns (H-N) across multiple sheets
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 28
$ws.Range("H28").Value = 1597.3636
$ws.Range("I28").Value = 1597.3636
$ws.Range("K28").Value = 1597.3636
$ws.Range("M28").Value = -1112.3636
# Row 40
$ws.Range("H40").Value = 2480.0476
$ws.Range("I40").Value = 2332.4167
$ws.Range("K40").Value = 2332.4167
$ws.Range("M40").Value = -2157.4167
# Row 112
$ws.Range("H112").Value = 1864.2858
$ws.Range("I112").Value = 1006.5
$ws.Range("J112").Value = 2207.4
$ws.Range("K112").Value = 3019.5
$ws.Range("L112").Value = 6622.200000000001
$ws.Range("M112").Value = -1911.5
$ws.Range("N112").Value = -8838.200000000001
# Row 113
$ws.Range("H113").Value = 3744
$ws.Range("I113").Value = 3528.2856
$ws.Range("K113").Value = 3528.2856
$ws.Range("M113").Value = -274.2856000000002
# Row 131
$ws.Range("H131").Value = 9619.615
$ws.Range("J131").Value = 19583.334
$ws.Range("L131").Value = 58750.00199999999
$ws.Range("N131").Value = -68830.00199999999
# Row 138
$ws.Range("H138").Value = 3474.87
$ws.Range("J138").Value = 3462.7705
$ws.Range("L138").Value = 10388.3115
$ws.Range("N138").Value = -20668.3115

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 5017.028
$ws.Range("I32").Value = 4775.8184
$ws.Range("J32").Value = 7670.3335
$ws.Range("K32").Value = 4775.8184
$ws.Range("L32").Value = 7670.3335
$ws.Range("M32").Value = -4488.8184
$ws.Range("N32").Value = -8244.333500000001
# Row 92
$ws.Range("H92").Value = 60448.5
$ws.Range("J92").Value = 60448.5
$ws.Range("L92").Value = 60448.5
$ws.Range("N92").Value = -65440.5
# Row 101
$ws.Range("H101").Value = 44989
$ws.Range("J101").Value = 44989
$ws.Range("L101").Value = 44989
$ws.Range("N101").Value = -51479
# Row 102
$ws.Range("H102").Value = 1499
$ws.Range("I102").Value = 1499
$ws.Range("K102").Value = 1499
$ws.Range("M102").Value = 123
# Row 103
$ws.Range("H103").Value = 105999.5
$ws.Range("J103").Value = 105999.5
$ws.Range("L103").Value = 105999.5
$ws.Range("N103").Value = -108343.5

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 64
$ws.Range("H64").Value = 434.14285
$ws.Range("I64").Value = 549
$ws.Range("J64").Value = 388.2
$ws.Range("K64").Value = 549
$ws.Range("L64").Value = 388.2
$ws.Range("M64").Value = -324
$ws.Range("N64").Value = -838.2
# Row 67
$ws.Range("H67").Value = 434.14285
$ws.Range("I67").Value = 549
$ws.Range("J67").Value = 388.2
$ws.Range("K67").Value = 549
$ws.Range("L67").Value = 388.2
$ws.Range("M67").Value = 231
$ws.Range("N67").Value = -1948.2
# Row 99
$ws.Range("H99").Value = 2836.1428
$ws.Range("I99").Value = 1870.8
$ws.Range("J99").Value = 5249.5
$ws.Range("K99").Value = 1870.8
$ws.Range("L99").Value = 5249.5
$ws.Range("M99").Value = -372.8
$ws.Range("N99").Value = -8245.5
# Row 107
$ws.Range("H107").Value = 2768.9473
$ws.Range("J107").Value = 4481
$ws.Range("L107").Value = 4481
$ws.Range("N107").Value = -8321
# Row 126
$ws.Range("H126").Value = 95999
$ws.Range("J126").Value = 95999
$ws.Range("L126").Value = 95999
$ws.Range("N126").Value = -105879
# Row 134
$ws.Range("H134").Value = 485.75
$ws.Range("I134").Value = 485.75
$ws.Range("K134").Value = 1457.25
$ws.Range("M134").Value = 1077.75

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 3382.25
$ws.Range("I31").Value = 3154.2856
$ws.Range("J31").Value = 3559.5557
$ws.Range("K31").Value = 3154.2856
$ws.Range("L31").Value = 3559.5557
$ws.Range("M31").Value = -2859.2856
$ws.Range("N31").Value = -4149.5557
# Row 34
$ws.Range("H34").Value = 3382.25
$ws.Range("I34").Value = 3154.2856
$ws.Range("J34").Value = 3559.5557
$ws.Range("K34").Value = 3154.2856
$ws.Range("L34").Value = 3559.5557
$ws.Range("M34").Value = -2952.2856
$ws.Range("N34").Value = -3963.5557
# Row 122
$ws.Range("H122").Value = 1292.3334
$ws.Range("I122").Value = 1292.3334
$ws.Range("K122").Value = 3877.0002
$ws.Range("M122").Value = -1427.0002
# Row 132
$ws.Range("H132").Value = 1520.6154
$ws.Range("I132").Value = 1253.1
$ws.Range("K132").Value = 3759.3
$ws.Range("M132").Value = -1229.3
# Row 134
$ws.Range("H134").Value = 1898.9656
$ws.Range("I134").Value = 1831.1428
$ws.Range("K134").Value = 5493.428400000001
$ws.Range("M134").Value = -2958.428400000001

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 21
$ws.Range("H21").Value = 95
$ws.Range("J21").Value = 90
$ws.Range("L21").Value = 270
$ws.Range("N21").Value = -616
# Row 82
$ws.Range("H82").Value = 2006.5
$ws.Range("I82").Value = 3013
$ws.Range("K82").Value = 9039
$ws.Range("M82").Value = -8633
# Row 85
$ws.Range("H85").Value = 2006.5
$ws.Range("I85").Value = 3013
$ws.Range("K85").Value = 9039
$ws.Range("M85").Value = -7635
# Row 92
$ws.Range("H92").Value = 799.3333
$ws.Range("I92").Value = 759.4
$ws.Range("J92").Value = 999
$ws.Range("K92").Value = 2278.2
$ws.Range("L92").Value = 2997
$ws.Range("M92").Value = -1030.2
$ws.Range("N92").Value = -5493
# Row 141
$ws.Range("H141").Value = 13624.5
$ws.Range("I141").Value = 13624.5
$ws.Range("K141").Value = 40873.5
$ws.Range("M141").Value = -35693.5

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 62
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").ClearContents()
# Row 65
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").ClearContents()
# Row 80
$ws.Range("H80").Value = 10247.375
$ws.Range("I80").Value = 2788.75
$ws.Range("J80").Value = 12733.583
$ws.Range("K80").Value = 2788.75
$ws.Range("L80").Value = 12733.583
$ws.Range("M80").Value = -1790.75
$ws.Range("N80").Value = -14729.583
# Row 83
$ws.Range("H83").Value = 10247.375
$ws.Range("I83").Value = 2788.75
$ws.Range("J83").Value = 12733.583
$ws.Range("K83").Value = 13943.75
$ws.Range("L83").Value = 63667.915
$ws.Range("M83").Value = -8951.75
$ws.Range("N83").Value = -73651.91500000001
# Row 104
$ws.Range("H104").Value = 49989
$ws.Range("J104").Value = 49989
$ws.Range("L104").Value = 49989
$ws.Range("N104").Value = -56977
# Row 113
$ws.Range("H113").Value = 1542
$ws.Range("I113").Value = 1542
$ws.Range("K113").Value = 1542
$ws.Range("M113").Value = 628
# Row 132
$ws.Range("H132").Value = 2999.5
$ws.Range("I132").Value = 2999.5
$ws.Range("K132").Value = 8998.5
$ws.Range("M132").Value = -6468.5

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 14
$ws.Range("H14").Value = 0
$ws.Range("J14").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("N14").ClearContents()
# Row 16
$ws.Range("H16").Value = 429.625
$ws.Range("I16").Value = 456.33334
$ws.Range("K16").Value = 456.33334
$ws.Range("M16").Value = -286.33334
# Row 43
$ws.Range("H43").Value = 49995
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 49995
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 49995
$ws.Range("M43").ClearContents()
$ws.Range("N43").Value = -50381
# Row 61
$ws.Range("H61").Value = 3571.4375
$ws.Range("I61").Value = 3422.12
$ws.Range("J61").Value = 4104.7144
$ws.Range("K61").Value = 3422.12
$ws.Range("L61").Value = 4104.7144
$ws.Range("M61").Value = -3220.12
$ws.Range("N61").Value = -4508.7144
# Row 68
$ws.Range("H68").Value = 2066.6667
$ws.Range("J68").Value = 3533.3333
$ws.Range("L68").Value = 3533.3333
$ws.Range("N68").Value = -5031.3333
# Row 71
$ws.Range("H71").Value = 2066.6667
$ws.Range("J71").Value = 3533.3333
$ws.Range("L71").Value = 17666.6665
$ws.Range("N71").Value = -25154.6665
# Row 100
$ws.Range("H100").Value = 4199.857
$ws.Range("I100").Value = 1249.5
$ws.Range("K100").Value = 1249.5
$ws.Range("M100").Value = -708.5
# Row 113
$ws.Range("H113").Value = 3571.4375
$ws.Range("I113").Value = 3422.12
$ws.Range("J113").Value = 4104.7144
$ws.Range("K113").Value = 3422.12
$ws.Range("L113").Value = 4104.7144
$ws.Range("M113").Value = -1252.12
$ws.Range("N113").Value = -8444.714400000001

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 9
$ws.Range("H9").Value = 7832.3335
$ws.Range("I9").Value = 10495
$ws.Range("J9").Value = 2507
$ws.Range("K9").Value = 10495
$ws.Range("L9").Value = 2507
$ws.Range("M9").Value = -10355
$ws.Range("N9").Value = -2787
# Row 14
$ws.Range("H14").Value = 966.55554
$ws.Range("I14").Value = 983.3333
$ws.Range("J14").Value = 933
$ws.Range("K14").Value = 983.3333
$ws.Range("L14").Value = 933
$ws.Range("M14").Value = -815.3333
$ws.Range("N14").Value = -1269
# Row 34
$ws.Range("H34").Value = 16681
$ws.Range("I34").Value = 13363
$ws.Range("J34").Value = 19999
$ws.Range("K34").Value = 13363
$ws.Range("L34").Value = 19999
$ws.Range("M34").Value = -13160
$ws.Range("N34").Value = -20405
# Row 70
$ws.Range("H70").Value = 0
$ws.Range("I70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("M70").ClearContents()
# Row 73
$ws.Range("H73").Value = 0
$ws.Range("I73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("M73").ClearContents()
# Row 75
$ws.Range("H75").Value = 80207.2
$ws.Range("I75").Value = 80207.2
$ws.Range("K75").Value = 80207.2
$ws.Range("M75").Value = -79271.2
# Row 78
$ws.Range("H78").Value = 80207.2
$ws.Range("I78").Value = 80207.2
$ws.Range("K78").Value = 240621.6
$ws.Range("M78").Value = -235941.6
# Row 100
$ws.Range("H100").Value = 20001556
$ws.Range("I100").Value = 33334344
$ws.Range("J100").Value = 2375
$ws.Range("K100").Value = 66668688
$ws.Range("L100").Value = 4750
$ws.Range("M100").Value = -66668147
$ws.Range("N100").Value = -5832
# Row 126
$ws.Range("H126").Value = 4175.0527
$ws.Range("I126").Value = 2368
$ws.Range("K126").Value = 7104
$ws.Range("M126").Value = -4634
